# "added header=False on excel output. Files upload to A2"
# The exporter no longer writes a "Key"/"Value" header row - the
# key/value pairs now start right at row 1 (previously row 2).
# Remove the old header row; Excel shifts all remaining rows (and their
# formatting) up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(1).Delete()
